$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2253968253968254
$ws.Range("C2").Value = 0.5047619047619047
$ws.Range("J2").Value = 0.0253968253968254
$ws.Range("P2").Value = 0.1523809523809524
$ws.Range("S2").Value = 0.09206349206349207
$ws.Range("B3").Value = 0.005952380952380952
$ws.Range("C3").Value = 0.03571428571428571
$ws.Range("J3").Value = 0.0119047619047619
$ws.Range("P3").Value = 0.7619047619047619
$ws.Range("S3").Value = 0.1845238095238095
$ws.Range("J4").Value = 0.103448275862069
$ws.Range("P4").Value = 0.7241379310344828
$ws.Range("S4").Value = 0.1724137931034483
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.75
$ws.Range("B6").Value = 0.07582938388625593
$ws.Range("D6").Value = 0.01895734597156398
$ws.Range("F6").Value = 0.06635071090047394
$ws.Range("J6").Value = 0.1658767772511848
$ws.Range("O6").Value = 0.04739336492890995
$ws.Range("Q6").Value = 0.1137440758293839
$ws.Range("R6").Value = 0.08530805687203792
$ws.Range("S6").Value = 0.4265402843601896
$ws.Range("B7").Value = 0.1261261261261261
$ws.Range("D7").Value = 0.01801801801801802
$ws.Range("E7").Value = 0.009009009009009009
$ws.Range("F7").Value = 0.02702702702702703
$ws.Range("J7").Value = 0.1351351351351351
$ws.Range("O7").Value = 0.01801801801801802
$ws.Range("Q7").Value = 0.1846846846846847
$ws.Range("R7").Value = 0.06756756756756757
$ws.Range("S7").Value = 0.4144144144144144
$ws.Range("B8").Value = 0.0989247311827957
$ws.Range("D8").Value = 0.01075268817204301
$ws.Range("E8").Value = 0.002150537634408602
$ws.Range("F8").Value = 0.04731182795698925
$ws.Range("J8").Value = 0.1032258064516129
$ws.Range("O8").Value = 0.02365591397849462
$ws.Range("Q8").Value = 0.178494623655914
$ws.Range("R8").Value = 0.07956989247311828
$ws.Range("S8").Value = 0.4559139784946237
$ws.Range("B9").Value = 0.07936507936507936
$ws.Range("D9").Value = 0.005291005291005291
$ws.Range("F9").Value = 0.1005291005291005
$ws.Range("J9").Value = 0.1481481481481481
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.164021164021164
$ws.Range("R9").Value = 0.06878306878306878
$ws.Range("S9").Value = 0.4232804232804233
$ws.Range("B10").Value = 0.1049334377447142
$ws.Range("D10").Value = 0.01252936570086139
$ws.Range("E10").Value = 0.001566170712607674
$ws.Range("F10").Value = 0.07909162098668755
$ws.Range("J10").Value = 0.1362568519968677
$ws.Range("O10").Value = 0.02427564604541895
$ws.Range("Q10").Value = 0.1918559122944401
$ws.Range("R10").Value = 0.08692247454972592
$ws.Range("S10").Value = 0.3625685199686766
$ws.Range("G11").Value = 0.138121546961326
$ws.Range("J11").Value = 0.1022099447513812
$ws.Range("K11").Value = 0.2044198895027624
$ws.Range("L11").Value = 0.5248618784530387
$ws.Range("S11").Value = 0.03038674033149171
$ws.Range("G12").Value = 0.7295918367346939
$ws.Range("J12").Value = 0.2193877551020408
$ws.Range("K12").Value = 0.01530612244897959
$ws.Range("L12").Value = 0.01020408163265306
$ws.Range("S12").Value = 0.02551020408163265
$ws.Range("G13").Value = 0.6964285714285714
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05357142857142857
$ws.Range("F15").Value = 0.02145922746781116
$ws.Range("H15").Value = 0.2103004291845494
$ws.Range("I15").Value = 0.06437768240343347
$ws.Range("J15").Value = 0.2746781115879828
$ws.Range("K15").Value = 0.05150214592274678
$ws.Range("M15").Value = 0.004291845493562232
$ws.Range("O15").Value = 0.06866952789699571
$ws.Range("S15").Value = 0.3047210300429185
$ws.Range("F16").Value = 0.01025641025641026
$ws.Range("H16").Value = 0.1692307692307692
$ws.Range("I16").Value = 0.09230769230769231
$ws.Range("J16").Value = 0.3692307692307693
$ws.Range("K16").Value = 0.1538461538461539
$ws.Range("M16").Value = 0.03076923076923077
$ws.Range("O16").Value = 0.04102564102564103
$ws.Range("S16").Value = 0.1333333333333333
$ws.Range("F17").Value = 0.007092198581560284
$ws.Range("H17").Value = 0.1607565011820331
$ws.Range("I17").Value = 0.1016548463356974
$ws.Range("J17").Value = 0.3995271867612293
$ws.Range("K17").Value = 0.1347517730496454
$ws.Range("M17").Value = 0.02364066193853428
$ws.Range("O17").Value = 0.05673758865248227
$ws.Range("S17").Value = 0.115839243498818
$ws.Range("F18").Value = 0.005208333333333333
$ws.Range("H18").Value = 0.2552083333333333
$ws.Range("I18").Value = 0.0625
$ws.Range("J18").Value = 0.4114583333333333
$ws.Range("K18").Value = 0.078125
$ws.Range("M18").Value = 0.01041666666666667
$ws.Range("O18").Value = 0.046875
$ws.Range("S18").Value = 0.1302083333333333
$ws.Range("F19").Value = 0.008245877061469266
$ws.Range("H19").Value = 0.2061469265367316
$ws.Range("I19").Value = 0.07721139430284858
$ws.Range("J19").Value = 0.3628185907046477
$ws.Range("K19").Value = 0.1259370314842579
$ws.Range("M19").Value = 0.02848575712143928
$ws.Range("N19").Value = 0.0007496251874062968
$ws.Range("O19").Value = 0.0704647676161919
$ws.Range("S19").Value = 0.1199400299850075
